$p = $ppt.ActivePresentation
$p.ApplyTheme("ppt/theme/theme2.xml")
Write-Output "done"
